$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the standalone "Meta description: ..." paragraph that sits
# right after the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: at the end of the document, insert a new bold paragraph
# ("Play Age of the Gods Free | Review of Greek Mythology Slot") right before
# the final (italic) paragraph, and replace the final paragraph's text (the
# old image-generation prompt) with the meta-description sentence, keeping
# its italic formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs($count - 1)
$r = $secondLast.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara = $d.Paragraphs($count)
$newRange = $newPara.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of the Gods Free | Review of Greek Mythology Slot</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newRange.InsertXML($xml)

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$oldText = 'Create a feature image for the game "Age of the Gods" that highlights the Greek mythology theme and features a happy Maya warrior with glasses. The image should be in cartoon style and should include Mount Olympus and the main deities such as Athena, Zeus, Poseidon, Hades, and Hercules. The Maya warrior can be positioned in the center of the image, holding a slot machine lever or spinning a wheel with a confident and excited expression. The overall tone should be vibrant and colorful, highlighting the adventurous and rewarding nature of this popular online slot game.'
$newText = 'Try Age of the Gods slot game for free! Our review covers gameplay features, design, progressive jackpots, accessibility, and more.'
$null = $lastRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
